$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update row 13 (currently holds placeholder "XXXXXXX" values) with the new
# Samoa Projection Data entry, mirroring row 12's lat/long/elevation values.
$ws.Range("B13").Value = "Samoa Projection Data"
$ws.Range("C13").Value = "Samoa_Climate_Projection.xlsx"
$ws.Range("D13").Value = -13.759
$ws.Range("E13").Value = -172.1046
$ws.Range("F13").Value = 23

# Move the active selection as recorded after the edit.
$ws.Activate()
$ws.Range("D28").Select()
